$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row labels: "_old" -> "_FV2310" and "_new" -> "_FV2404"
$headerMap = @{
    "A1" = "Segmentname_FV2310"
    "B1" = "Segmentgruppe_FV2310"
    "C1" = "Segment_FV2310"
    "D1" = "Datenelement_FV2310"
    "E1" = "Segment ID_FV2310"
    "F1" = "Code_FV2310"
    "G1" = "Qualifier_FV2310"
    "H1" = "Beschreibung_FV2310"
    "I1" = "Bedingungsausdruck_FV2310"
    "J1" = "Bedingung_FV2310"
    "L1" = "Segmentname_FV2404"
    "M1" = "Segmentgruppe_FV2404"
    "N1" = "Segment_FV2404"
    "O1" = "Datenelement_FV2404"
    "P1" = "Segment ID_FV2404"
    "Q1" = "Code_FV2404"
    "R1" = "Qualifier_FV2404"
    "S1" = "Beschreibung_FV2404"
    "T1" = "Bedingungsausdruck_FV2404"
    "U1" = "Bedingung_FV2404"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# Turn the used range into an Excel table (ListObject)
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U67"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
